$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "DESCRIPCIÓN" column is inserted right after "CÓDIGO PROVEEDOR" (i.e. before the
# old column B "CANTIDAD"), shifting CANTIDAD/PRECIO/FOLIO/FECHA EN LA FACTURA one column right.
$ws.Columns("B:B").Insert()

# Header text for the freshly inserted column (adds a new shared string).
$ws.Range("B1").Value = "DESCRIPCIÓN"

# Match the new column's width (as close as this runtime's column-width rounding allows).
$ws.Columns("B:B").ColumnWidth = 36.333333333333336

# The active selection moved from F5 to B5.
$ws.Range("B5").Select()
